# Fruta / hortaliza, semanal
# Insert a new weekly data row at sheet row 320 (pushes existing rows 320:343
# down to 321:344) and populate the new row with this week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 320:343 down to 321:344, leaving a blank row 320 in place.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new weekly record.
$ws.Range("A320").Value = 3
$ws.Range("B320").Value = "Femacal de La Calera"
$ws.Range("C320").Value = "Coquimbo"
$ws.Range("D320").Value = 44585
$ws.Range("E320").Value = 5
$ws.Range("F320").Value = 100112032
$ws.Range("G320").Value = "Zapallo italiano"
$ws.Range("H320").Value = "Sin especificar"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 165
$ws.Range("K320").Value = 9000
$ws.Range("L320").Value = 9500
$ws.Range("M320").Value = 9348
$ws.Range("N320").Value = "`$/caja 70 unidades"
$ws.Range("O320").Value = "Limache"
$ws.Range("P320").Value = 134
$ws.Range("Q320").Value = 70
$ws.Range("R320").Value = "Hortaliza"
